$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Totaal" row (row 2) with new aggregated values
$ws.Range("B2").Value = 919
$ws.Range("C2").Value = 9160
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 9190
$ws.Range("F2").Value = 2969
$ws.Range("G2").Value = 5223
$ws.Range("H2").Value = 0.01
$ws.Range("I2").Value = $false

# Add a new row 4 with the "test lokaal" data
$ws.Range("A4").Value = "test lokaal"
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 80
$ws.Range("D4").Value = 30
$ws.Range("E4").Value = 110
$ws.Range("F4").Value = 17
$ws.Range("G4").Value = 83
$ws.Range("H4").Value = 0.05
$ws.Range("I4").Value = $true

# Keep the "numbers stored as text" error-checking suppression in sync with
# the table's new extent (A1:I3 -> A1:I4).
$ws.Range("A1:I4").Errors.Item(9).Ignore = $true
